# Auto-generated cell updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 98524.62
$ws.Range("I15").Value = 98524.62
$ws.Range("K15").Value = 295573.86
$ws.Range("M15").Value = -295404.86
$ws.Range("H74").Value = 3423
$ws.Range("I74").Value = 2793.3333
$ws.Range("J74").Value = 3800.8
$ws.Range("K74").Value = 2793.3333
$ws.Range("L74").Value = 3800.8
$ws.Range("M74").Value = -1857.3333
$ws.Range("N74").Value = -5672.8
$ws.Range("H77").Value = 3423
$ws.Range("I77").Value = 2793.3333
$ws.Range("J77").Value = 3800.8
$ws.Range("K77").Value = 13966.6665
$ws.Range("L77").Value = 19004
$ws.Range("M77").Value = -9286.666499999999
$ws.Range("N77").Value = -28364
$ws.Range("H98").Value = 1023911.4
$ws.Range("I98").Value = 1125480.9
$ws.Range("J98").Value = 8216
$ws.Range("K98").Value = 1125480.9
$ws.Range("L98").Value = 8216
$ws.Range("M98").Value = -1123982.9
$ws.Range("N98").Value = -11212
$ws.Range("H122").Value = 1023911.4
$ws.Range("I122").Value = 1125480.9
$ws.Range("J122").Value = 8216
$ws.Range("K122").Value = 3376442.7
$ws.Range("L122").Value = 24648
$ws.Range("M122").Value = -3373992.7
$ws.Range("N122").Value = -29548
$ws.Range("H125").Value = 16017492
$ws.Range("I125").Value = 1333
$ws.Range("J125").Value = 37372372
$ws.Range("K125").Value = 11997
$ws.Range("L125").Value = 336351348
$ws.Range("M125").Value = -9537
$ws.Range("N125").Value = -336356268
$ws.Range("H141").Value = 4587.4
$ws.Range("I141").Value = 2971.5625
$ws.Range("J141").Value = 7460
$ws.Range("K141").Value = 8914.6875
$ws.Range("L141").Value = 22380
$ws.Range("M141").Value = -3734.6875
$ws.Range("N141").Value = -32740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1029.4615
$ws.Range("I4").Value = 1298.2222
$ws.Range("K4").Value = 1298.2222
$ws.Range("M4").Value = -1182.2222
$ws.Range("H5").Value = 250903
$ws.Range("I5").Value = 334403.34
$ws.Range("J5").Value = 402
$ws.Range("K5").Value = 334403.34
$ws.Range("L5").Value = 402
$ws.Range("M5").Value = -334291.34
$ws.Range("N5").Value = -626
$ws.Range("H6").Value = 5334
$ws.Range("I6").Value = 10002
$ws.Range("K6").Value = 10002
$ws.Range("M6").Value = -9829
$ws.Range("H32").Value = 23610.322
$ws.Range("I32").Value = 5487.2593
$ws.Range("J32").Value = 145941
$ws.Range("K32").Value = 5487.2593
$ws.Range("L32").Value = 145941
$ws.Range("M32").Value = -5200.2593
$ws.Range("N32").Value = -146515
$ws.Range("H109").Value = 30911
$ws.Range("J109").Value = 30911
$ws.Range("L109").Value = 30911
$ws.Range("N109").Value = -33685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250903
$ws.Range("I4").Value = 334403.34
$ws.Range("J4").Value = 402
$ws.Range("K4").Value = 334403.34
$ws.Range("L4").Value = 402
$ws.Range("M4").Value = -334288.34
$ws.Range("N4").Value = -632
$ws.Range("H15").Value = 4000
$ws.Range("J15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("N15").Value = -4454
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H94").Value = 655.46875
$ws.Range("I94").Value = 508.36
$ws.Range("J94").Value = 1180.8572
$ws.Range("K94").Value = 508.36
$ws.Range("L94").Value = 1180.8572
$ws.Range("M94").Value = -57.36000000000001
$ws.Range("N94").Value = -2082.8572
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3233.5557
$ws.Range("I31").Value = 1389.6
$ws.Range("J31").Value = 6630.316
$ws.Range("K31").Value = 1389.6
$ws.Range("L31").Value = 6630.316
$ws.Range("M31").Value = -1094.6
$ws.Range("N31").Value = -7220.316
$ws.Range("H34").Value = 3233.5557
$ws.Range("I34").Value = 1389.6
$ws.Range("J34").Value = 6630.316
$ws.Range("K34").Value = 1389.6
$ws.Range("L34").Value = 6630.316
$ws.Range("M34").Value = -1187.6
$ws.Range("N34").Value = -7034.316
$ws.Range("H58").Value = 1710.7858
$ws.Range("I58").Value = 1185.2174
$ws.Range("J58").Value = 4128.4
$ws.Range("K58").Value = 1185.2174
$ws.Range("L58").Value = 4128.4
$ws.Range("M58").Value = -982.2174
$ws.Range("N58").Value = -4534.4
$ws.Range("H136").Value = 1710.7858
$ws.Range("I136").Value = 1185.2174
$ws.Range("J136").Value = 4128.4
$ws.Range("K136").Value = 3555.6522
$ws.Range("L136").Value = 12385.2
$ws.Range("M136").Value = -1005.6522
$ws.Range("N136").Value = -17485.2
$ws.Range("H141").Value = 207274.97
$ws.Range("J141").Value = 207274.97
$ws.Range("L141").Value = 207274.97
$ws.Range("N141").Value = -217634.97

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 473.4091
$ws.Range("J107").Value = 493.4
$ws.Range("L107").Value = 1480.2
$ws.Range("N107").Value = -5320.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 28000
$ws.Range("J93").Value = 28000
$ws.Range("L93").Value = 28000
$ws.Range("N93").Value = -31744
$ws.Range("H123").Value = 11322.5
$ws.Range("J123").Value = 11322.5
$ws.Range("L123").Value = 11322.5
$ws.Range("N123").Value = -16222.5
$ws.Range("H132").Value = 4327.3706
$ws.Range("I132").Value = 4178.2856
$ws.Range("J132").Value = 4487.923
$ws.Range("K132").Value = 12534.8568
$ws.Range("L132").Value = 13463.769
$ws.Range("M132").Value = -10004.8568
$ws.Range("N132").Value = -18523.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17333.166
$ws.Range("I22").Value = 499.75
$ws.Range("J22").Value = 51000
$ws.Range("K22").Value = 499.75
$ws.Range("L22").Value = 51000
$ws.Range("M22").Value = -204.75
$ws.Range("N22").Value = -51590
$ws.Range("H27").Value = 17333.166
$ws.Range("I27").Value = 499.75
$ws.Range("J27").Value = 51000
$ws.Range("K27").Value = 499.75
$ws.Range("L27").Value = 51000
$ws.Range("M27").Value = -392.75
$ws.Range("N27").Value = -51214
$ws.Range("H55").Value = 292.5
$ws.Range("I55").Value = 282.2857
$ws.Range("J55").Value = 328.25
$ws.Range("K55").Value = 282.2857
$ws.Range("L55").Value = 328.25
$ws.Range("M55").Value = -109.2857
$ws.Range("N55").Value = -674.25
$ws.Range("H100").Value = 2302.8293
$ws.Range("I100").Value = 1494.1428
$ws.Range("J100").Value = 2722.1482
$ws.Range("K100").Value = 1494.1428
$ws.Range("L100").Value = 2722.1482
$ws.Range("M100").Value = -953.1428000000001
$ws.Range("N100").Value = -3804.1482
$ws.Range("H122").Value = 3696.138
$ws.Range("I122").Value = 3069.7144
$ws.Range("J122").Value = 3895.4546
$ws.Range("K122").Value = 9209.143199999999
$ws.Range("L122").Value = 11686.3638
$ws.Range("M122").Value = -6759.143199999999
$ws.Range("N122").Value = -16586.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3275.7568
$ws.Range("I81").Value = 2171.4707
$ws.Range("J81").Value = 4214.4
$ws.Range("K81").Value = 4342.9414
$ws.Range("L81").Value = 8428.799999999999
$ws.Range("M81").Value = -3281.9414
$ws.Range("N81").Value = -10550.8
$ws.Range("H84").Value = 3275.7568
$ws.Range("I84").Value = 2171.4707
$ws.Range("J84").Value = 4214.4
$ws.Range("K84").Value = 21714.707
$ws.Range("L84").Value = 42144
$ws.Range("M84").Value = -16410.707
$ws.Range("N84").Value = -52752

# Cells removed entirely in the update (no longer present in the row)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N92").ClearContents()
$ws.Range("N95").ClearContents()
